# Insert a new price-record row for "Betarraga" at Excel row 612 (pushing the
# existing row 612 and everything below it down by one row: 612->613, ...,
# 672->673). The newly inserted row 612 duplicates the data that used to live
# in row 612, except for the "Fecha" (D) and "Volumen" (J) values, which get
# the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastCol = 18   # columns A..R
$insertRow = 612

# Snapshot the row that is about to be pushed down - this becomes the
# template for the brand-new row we are inserting.
$template = @()
for ($c = 1; $c -le $lastCol; $c++) {
    $template += , $ws.Cells.Item($insertRow, $c).Value2
}

# Push row 612 (and everything after it) down by one row.
$ws.Rows.Item($insertRow).Insert()

# Populate the freshly inserted (now blank) row 612 with the template values.
for ($c = 1; $c -le $lastCol; $c++) {
    $ws.Cells.Item($insertRow, $c).Value = $template[$c - 1]
}

# Overwrite the two fields that differ for this new record:
#   D = Fecha (date serial), J = Volumen
$ws.Cells.Item($insertRow, 4).Value = 45194
$ws.Cells.Item($insertRow, 10).Value = 30
